# Updated symbol list on Tue Jan 31 19:19:05 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '312.86'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.44%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '37.62'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.88%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.134'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.22%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07914'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.86%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.906'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-2.22%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.269'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.48%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.850'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-7.58%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9217'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.19%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1225'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-3.57%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1922'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.38%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09134'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '5.08%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03310'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-4.56%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09627'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.58%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001386'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.96%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005709'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-5.55%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.514'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-2.21%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.416'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.19%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.45%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.277'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '4.50%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1272'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.13%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2589'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.04%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.84%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04368'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.15%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.56%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004311'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-3.94%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001221'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-10.43%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02229'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2.15%'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.91%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007391'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-3.17%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1361'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2.03%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008750'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-11.16%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002012'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.10%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008648'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.96%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006725'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-2.51%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.85%'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003362'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '10.75%'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001199'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-8.60%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002102'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.85%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002002'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.85%'
